$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.95417514485077
$ws.Range("C2").Value = 11.38503740958157
$ws.Range("E2").Value = 16.60989560622437
$ws.Range("F2").Value = 16.86991607391245
$ws.Range("G2").Value = 3.584089804427038
$ws.Range("I2").Value = 15.07788214228017
$ws.Range("O2").Value = 15.62271078555588

$ws.Range("B3").Value = 14.09972746306399
$ws.Range("C3").Value = 10.74138352244371
$ws.Range("E3").Value = 15.65977732696368
$ws.Range("F3").Value = 15.89584955866815
$ws.Range("G3").Value = 3.586345421855879
$ws.Range("I3").Value = 15.27051264362134
$ws.Range("O3").Value = 15.76623067272593

$ws.Range("B4").Value = 13.54700979990085
$ws.Range("C4").Value = 10.32420638013924
$ws.Range("E4").Value = 15.05097772297288
$ws.Range("F4").Value = 15.26997757108489
$ws.Range("G4").Value = 3.587799831309748
$ws.Range("I4").Value = 15.39456118671229
$ws.Range("O4").Value = 15.8617997034256

$ws.Range("B5").Value = 13.3148480874249
$ws.Range("C5").Value = 10.14875461013006
$ws.Range("E5").Value = 14.79676021549598
$ws.Range("F5").Value = 15.008197319934
$ws.Range("G5").Value = 3.588410041314026
$ws.Range("I5").Value = 15.44656712586401
$ws.Range("O5").Value = 15.90260070925761

$ws.Range("B6").Value = 13.2758841948271
$ws.Range("C6").Value = 10.11929447645399
$ws.Range("E6").Value = 14.75418638070411
$ws.Range("F6").Value = 14.96433081551589
$ws.Range("G6").Value = 3.588512426596147
$ws.Range("I6").Value = 15.45529066603848
$ws.Range("O6").Value = 15.90948728872328

$ws.Range("B7").Value = 13.54390662024748
$ws.Range("C7").Value = 10.32186212384934
$ws.Range("E7").Value = 15.04757367285683
$ws.Range("F7").Value = 15.26647399323133
$ws.Range("G7").Value = 3.587807989769907
$ws.Range("I7").Value = 15.3952566604473
$ws.Range("O7").Value = 15.86234246787765

$ws.Range("B8").Value = 14.6654632863702
$ws.Range("C8").Value = 11.16770589463535
$ws.Range("E8").Value = 16.28771535500158
$ws.Range("F8").Value = 16.5399640634477
$ws.Range("G8").Value = 3.584853164314834
$ws.Range("I8").Value = 15.14310517619997
$ws.Range("O8").Value = 15.67064055568737

$ws.Range("B9").Value = 16.63739350211652
$ws.Range("C9").Value = 12.64956017594646
$ws.Range("E9").Value = 18.63158843625158
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 3.57960696882296
$ws.Range("I9").Value = 14.69427648299733
$ws.Range("O9").Value = 15.35457154223684

$ws.Range("B10").Value = 17.94292611747911
$ws.Range("C10").Value = 13.62804050919103
$ws.Range("E10").Value = 20.29044273602282
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 3.576082771683159
$ws.Range("I10").Value = 14.39211588658997
$ws.Range("O10").Value = 15.15991968817342

$ws.Range("B11").Value = 18.50512727330458
$ws.Range("C11").Value = 14.04893390274913
$ws.Range("E11").Value = 21.00263295508407
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 3.574550359053559
$ws.Range("I11").Value = 14.26059960730766
$ws.Range("O11").Value = 15.07976187252135

$ws.Range("B12").Value = 18.71342606654253
$ws.Range("C12").Value = 14.20481596942529
$ws.Range("E12").Value = 21.26625698460389
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 3.573980184237847
$ws.Range("I12").Value = 14.21164840560022
$ws.Range("O12").Value = 15.05063490433985

$ws.Range("B13").Value = 18.66877002364902
$ws.Range("C13").Value = 14.1713999134305
$ws.Range("E13").Value = 21.20975011417035
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 3.574102532584222
$ws.Range("I13").Value = 14.22215310928772
$ws.Range("O13").Value = 15.05685301384862

$ws.Range("B14").Value = 18.52235639932158
$ws.Range("C14").Value = 14.06182866356835
$ws.Range("E14").Value = 21.02444285614956
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 3.574503248013778
$ws.Range("I14").Value = 14.25655531917401
$ws.Range("O14").Value = 15.07734086846188

$ws.Range("B15").Value = 18.43207454850534
$ws.Range("C15").Value = 13.99425672905252
$ws.Range("E15").Value = 20.91014773613494
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 3.574750013518217
$ws.Range("I15").Value = 14.27773843324543
$ws.Range("O15").Value = 15.09005068522104

$ws.Range("B16").Value = 17.90554379007511
$ws.Range("C16").Value = 13.60004504165616
$ws.Range("E16").Value = 20.24304857140823
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 3.576184337327184
$ws.Range("I16").Value = 14.40083003788647
$ws.Range("O16").Value = 15.16532877267049

$ws.Range("B17").Value = 17.57438709643368
$ws.Range("C17").Value = 13.35199021735866
$ws.Range("E17").Value = 19.82296192863376
$ws.Range("F17").Value = 20.20408069617459
$ws.Range("G17").Value = 3.577082330435755
$ws.Range("I17").Value = 14.47786165393236
$ws.Range("O17").Value = 15.21367290750535

$ws.Range("B18").Value = 17.38093297655628
$ws.Range("C18").Value = 13.20703636654496
$ws.Range("E18").Value = 19.57734334985946
$ws.Range("F18").Value = 19.95656407809808
$ws.Range("G18").Value = 3.577605496318012
$ws.Range("I18").Value = 14.52272727904176
$ws.Range("O18").Value = 15.24226765903491

$ws.Range("B19").Value = 17.31492196514286
$ws.Range("C19").Value = 13.15756666726852
$ws.Range("E19").Value = 19.49349377146694
$ws.Range("F19").Value = 19.87204792380562
$ws.Range("G19").Value = 3.577783777642642
$ws.Range("I19").Value = 14.5380141160499
$ws.Range("O19").Value = 15.25208415274904

$ws.Range("B20").Value = 17.60994817440042
$ws.Range("C20").Value = 13.37863212318749
$ws.Range("E20").Value = 19.86809417261692
$ws.Range("F20").Value = 20.24955283636157
$ws.Range("G20").Value = 3.576986048268385
$ws.Range("I20").Value = 14.46960366498036
$ws.Range("O20").Value = 15.20844483921919

$ws.Range("B21").Value = 18.56548657119015
$ws.Range("C21").Value = 14.09410752869517
$ws.Range("E21").Value = 21.07903647064006
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 3.574385274085989
$ws.Range("I21").Value = 14.24642747711277
$ws.Range("O21").Value = 15.07128962163519

$ws.Range("B22").Value = 19.16319505786531
$ws.Range("C22").Value = 14.54130052123812
$ws.Range("E22").Value = 21.83511032318818
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 3.57274445784505
$ws.Range("I22").Value = 14.10552897022255
$ws.Range("O22").Value = 14.98881273761024

$ws.Range("B23").Value = 18.84664789177977
$ws.Range("C23").Value = 14.30449729209777
$ws.Range("E23").Value = 21.43480236214343
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 3.573614818754262
$ws.Range("I23").Value = 14.18027617993392
$ws.Range("O23").Value = 15.03217001243691

$ws.Range("B24").Value = 17.59388056340996
$ws.Range("C24").Value = 13.36659461670327
$ws.Range("E24").Value = 19.84770267479086
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 3.577029555946285
$ws.Range("I24").Value = 14.47333529762623
$ws.Range("O24").Value = 15.21080595518124

$ws.Range("B25").Value = 16.12883122015734
$ws.Range("C25").Value = 12.26790570296586
$ws.Range("E25").Value = 17.98284373560391
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 3.580967930674202
$ws.Range("I25").Value = 14.81083458789333
$ws.Range("O25").Value = 15.4335590132769
